$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: A7's timestamp value is refreshed by the scheduled task run.
$ws.Range("A7").Value = 45873.58364391204

# New row 8: appended by the automatic scheduled update.
$ws.Range("A8").Value = 45873.62529006821
$ws.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B8").Value = 2025
$ws.Range("C8").Value = 15
$ws.Range("D8").Value = 20.56
$ws.Range("E8").Value = 74.34999999999999
$ws.Range("F8").Value = 493.95
$ws.Range("G8").Value = 14
$ws.Range("H8").Value = "ESE"
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = "15:00:25"
